# Scheduled-runner data refresh: updates currentAveragePrice* / LevePrice* /
# LeveProfit* columns (H:N) for a set of leve rows across all eight
# crafting-job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR), matching a
# fresh pull from the market data source. A few rows drop to/from a
# zero-NQ-price or zero-HQ-price state, which means the corresponding
# LeveProfitNQ (M) or LeveProfitHQ (N) cell must be cleared or created
# outright (the upstream generator omits that column's cell entirely when
# its source price is 0).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 13339.667  # was 22000
$ws.Range("I21").Value = 0  # was 50000
$ws.Range("J21").Value = 13339.667  # was 15000
$ws.Range("K21").Value = 0  # was 50000
$ws.Range("L21").Value = 13339.667  # was 15000
$ws.Range("M21").ClearContents()  # was -49532
$ws.Range("N21").Value = -14275.667  # was -15936

$ws.Range("H23").Value = 13339.667  # was 22000
$ws.Range("I23").Value = 0  # was 50000
$ws.Range("J23").Value = 13339.667  # was 15000
$ws.Range("K23").Value = 0  # was 50000
$ws.Range("L23").Value = 13339.667  # was 15000
$ws.Range("M23").ClearContents()  # was -49766
$ws.Range("N23").Value = -13807.667  # was -15468

$ws.Range("H137").Value = 41667984  # was 47620410
$ws.Range("I137").Value = 66667970  # was 83334710
$ws.Range("K137").Value = 200003910  # was 250004130
$ws.Range("M137").Value = -200001360  # was -250001580

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 25030.314  # was 29800.979
$ws.Range("I32").Value = 4480.6304  # was 5522.081
$ws.Range("J32").Value = 143191  # was 129614.22
$ws.Range("K32").Value = 4480.6304  # was 5522.081
$ws.Range("L32").Value = 143191  # was 129614.22
$ws.Range("M32").Value = -4193.6304  # was -5235.081
$ws.Range("N32").Value = -143765  # was -130188.22

$ws.Range("H61").Value = 3106.6885  # was 3078.6333
$ws.Range("I61").Value = 2389.9185  # was 2366.32
$ws.Range("J61").Value = 6033.5  # was 6640.2
$ws.Range("K61").Value = 2389.9185  # was 2366.32
$ws.Range("L61").Value = 6033.5  # was 6640.2
$ws.Range("M61").Value = -2177.9185  # was -2154.32
$ws.Range("N61").Value = -6457.5  # was -7064.2

$ws.Range("H74").Value = 8731.588  # was 7501.75
$ws.Range("I74").Value = 1316.3334  # was 1124.5
$ws.Range("J74").Value = 12776.272  # was 11753.25
$ws.Range("K74").Value = 1316.3334  # was 1124.5
$ws.Range("L74").Value = 12776.272  # was 11753.25
$ws.Range("M74").Value = -442.3334  # was -250.5
$ws.Range("N74").Value = -14524.272  # was -13501.25

$ws.Range("H77").Value = 8731.588  # was 7501.75
$ws.Range("I77").Value = 1316.3334  # was 1124.5
$ws.Range("J77").Value = 12776.272  # was 11753.25
$ws.Range("K77").Value = 6581.666999999999  # was 5622.5
$ws.Range("L77").Value = 63881.36  # was 58766.25
$ws.Range("M77").Value = -2213.666999999999  # was -1254.5
$ws.Range("N77").Value = -72617.36  # was -67502.25

$ws.Range("H132").Value = 2819.975  # was 3509.8667
$ws.Range("I132").Value = 2287.3794  # was 2977.6667
$ws.Range("J132").Value = 4224.091  # was 4751.6665
$ws.Range("K132").Value = 6862.138199999999  # was 8933.000100000001
$ws.Range("L132").Value = 12672.273  # was 14254.9995
$ws.Range("M132").Value = -4332.138199999999  # was -6403.000100000001
$ws.Range("N132").Value = -17732.273  # was -19314.9995

$ws.Range("H136").Value = 3106.6885  # was 3078.6333
$ws.Range("I136").Value = 2389.9185  # was 2366.32
$ws.Range("J136").Value = 6033.5  # was 6640.2
$ws.Range("K136").Value = 7169.755500000001  # was 7098.960000000001
$ws.Range("L136").Value = 18100.5  # was 19920.6
$ws.Range("M136").Value = -4619.755500000001  # was -4548.960000000001
$ws.Range("N136").Value = -23200.5  # was -25020.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 42332.332  # was 26992.334
$ws.Range("I26").Value = 28999  # was 26992.334
$ws.Range("J26").Value = 48999  # was 0
$ws.Range("K26").Value = 28999  # was 26992.334
$ws.Range("L26").Value = 48999  # was 0
$ws.Range("M26").Value = -28707  # was -26700.334
$ws.Range("N26").Value = -49583  # was None

$ws.Range("H134").Value = 3550  # was 2943.413
$ws.Range("I134").Value = 2446.4583  # was 1882.2941
$ws.Range("J134").Value = 5757.0835  # was 5949.9165
$ws.Range("K134").Value = 7339.374899999999  # was 5646.8823
$ws.Range("L134").Value = 17271.2505  # was 17849.7495
$ws.Range("M134").Value = -4804.374899999999  # was -3111.8823
$ws.Range("N134").Value = -22341.2505  # was -22919.7495

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4586.091  # was 4802.3076
$ws.Range("I31").Value = 1689.3889  # was 1893.375
$ws.Range("J31").Value = 5995.2974  # was 6095.1665
$ws.Range("K31").Value = 1689.3889  # was 1893.375
$ws.Range("L31").Value = 5995.2974  # was 6095.1665
$ws.Range("M31").Value = -1394.3889  # was -1598.375
$ws.Range("N31").Value = -6585.2974  # was -6685.1665

$ws.Range("H34").Value = 4586.091  # was 4802.3076
$ws.Range("I34").Value = 1689.3889  # was 1893.375
$ws.Range("J34").Value = 5995.2974  # was 6095.1665
$ws.Range("K34").Value = 1689.3889  # was 1893.375
$ws.Range("L34").Value = 5995.2974  # was 6095.1665
$ws.Range("M34").Value = -1487.3889  # was -1691.375
$ws.Range("N34").Value = -6399.2974  # was -6499.1665

$ws.Range("H58").Value = 3040.35  # was 2181.4243
$ws.Range("I58").Value = 1812.5  # was 1099.375
$ws.Range("J58").Value = 3858.9167  # was 3199.8235
$ws.Range("K58").Value = 1812.5  # was 1099.375
$ws.Range("L58").Value = 3858.9167  # was 3199.8235
$ws.Range("M58").Value = -1609.5  # was -896.375
$ws.Range("N58").Value = -4264.9167  # was -3605.8235

$ws.Range("H99").Value = 4816307  # was 7720
$ws.Range("I99").Value = 8941070  # was 17100
$ws.Range("J99").Value = 4083.3333  # was 3030
$ws.Range("K99").Value = 8941070  # was 17100
$ws.Range("L99").Value = 4083.3333  # was 3030
$ws.Range("M99").Value = -8939572  # was -15602
$ws.Range("N99").Value = -7079.3333  # was -6026

$ws.Range("H122").Value = 982.4761999999999  # was 987.4286
$ws.Range("I122").Value = 970.1111  # was 985.3333
$ws.Range("J122").Value = 1056.6666  # was 1000
$ws.Range("K122").Value = 2910.3333  # was 2955.9999
$ws.Range("L122").Value = 3169.9998  # was 3000
$ws.Range("M122").Value = -460.3332999999998  # was -505.9998999999998
$ws.Range("N122").Value = -8069.9998  # was -7900

$ws.Range("H126").Value = 4816307  # was 7720
$ws.Range("I126").Value = 8941070  # was 17100
$ws.Range("J126").Value = 4083.3333  # was 3030
$ws.Range("K126").Value = 26823210  # was 51300
$ws.Range("L126").Value = 12249.9999  # was 9090
$ws.Range("M126").Value = -26820740  # was -48830
$ws.Range("N126").Value = -17189.9999  # was -14030

$ws.Range("H134").Value = 4031.5264  # was 4859.357
$ws.Range("I134").Value = 1751.5  # was 2284.8
$ws.Range("J134").Value = 5689.727  # was 6289.6665
$ws.Range("K134").Value = 5254.5  # was 6854.400000000001
$ws.Range("L134").Value = 17069.181  # was 18868.9995
$ws.Range("M134").Value = -2719.5  # was -4319.400000000001
$ws.Range("N134").Value = -22139.181  # was -23938.9995

$ws.Range("H136").Value = 3040.35  # was 2181.4243
$ws.Range("I136").Value = 1812.5  # was 1099.375
$ws.Range("J136").Value = 3858.9167  # was 3199.8235
$ws.Range("K136").Value = 5437.5  # was 3298.125
$ws.Range("L136").Value = 11576.7501  # was 9599.470499999999
$ws.Range("M136").Value = -2887.5  # was -748.125
$ws.Range("N136").Value = -16676.7501  # was -14699.4705

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 6945574  # was 5748239
$ws.Range("I132").Value = 1002  # was 977
$ws.Range("J132").Value = 8334488.5  # was 7247524.5
$ws.Range("K132").Value = 9018  # was 8793
$ws.Range("L132").Value = 75010396.5  # was 65227720.5
$ws.Range("M132").Value = -6488  # was -6263
$ws.Range("N132").Value = -75015456.5  # was -65232780.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1676.7059  # was 1846.7142
$ws.Range("I122").Value = 1576.7693  # was 1713.4546
$ws.Range("J122").Value = 2001.5  # was 2335.3333
$ws.Range("K122").Value = 4730.3079  # was 5140.3638
$ws.Range("L122").Value = 6004.5  # was 7005.999899999999
$ws.Range("M122").Value = -2280.3079  # was -2690.3638
$ws.Range("N122").Value = -10904.5  # was -11905.9999

$ws.Range("H132").Value = 2604.84  # was 3708.0715
$ws.Range("I132").Value = 2748.8  # was 4302
$ws.Range("J132").Value = 2388.9  # was 2916.1667
$ws.Range("K132").Value = 8246.400000000001  # was 12906
$ws.Range("L132").Value = 7166.700000000001  # was 8748.500100000001
$ws.Range("M132").Value = -5716.400000000001  # was -10376
$ws.Range("N132").Value = -12226.7  # was -13808.5001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3153.8462  # was 3004
$ws.Range("I7").Value = 2411.111  # was 1971.4286
$ws.Range("J7").Value = 3547.0588  # was 3405.5557
$ws.Range("K7").Value = 2411.111  # was 1971.4286
$ws.Range("L7").Value = 3547.0588  # was 3405.5557
$ws.Range("M7").Value = -2299.111  # was -1859.4286
$ws.Range("N7").Value = -3771.0588  # was -3629.5557

$ws.Range("H126").Value = 3153.8462  # was 3004
$ws.Range("I126").Value = 2411.111  # was 1971.4286
$ws.Range("J126").Value = 3547.0588  # was 3405.5557
$ws.Range("K126").Value = 7233.333  # was 5914.2858
$ws.Range("L126").Value = 10641.1764  # was 10216.6671
$ws.Range("M126").Value = -4763.333  # was -3444.2858
$ws.Range("N126").Value = -15581.1764  # was -15156.6671

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 26499.5  # was 26500
$ws.Range("J12").Value = 2999  # was 3000
$ws.Range("L12").Value = 2999  # was 3000
$ws.Range("N12").Value = -3283  # was -3284

$ws.Range("H107").Value = 316.27777  # was 304.3684
$ws.Range("I107").Value = 245.61539  # was 234.5
$ws.Range("K107").Value = 736.84617  # was 703.5
$ws.Range("M107").Value = 1183.15383  # was 1216.5

$ws.Range("H126").Value = 48328.332  # was 46293.684
$ws.Range("I126").Value = 67040.336  # was 77351.62
$ws.Range("J126").Value = 1548.3334  # was 1432.2222
$ws.Range("K126").Value = 201121.008  # was 232054.86
$ws.Range("L126").Value = 4645.0002  # was 4296.6666
$ws.Range("M126").Value = -198651.008  # was -229584.86
$ws.Range("N126").Value = -9585.0002  # was -9236.6666

$ws.Range("H132").Value = 1797.2587  # was 2789.4182
$ws.Range("I132").Value = 1657.4565  # was 2916.7441
$ws.Range("K132").Value = 4972.3695  # was 8750.2323
$ws.Range("M132").Value = -2442.3695  # was -6220.2323

$ws.Range("H136").Value = 1526.0625  # was 1717.6786
$ws.Range("I136").Value = 719.2381  # was 845
$ws.Range("K136").Value = 2157.7143  # was 2535
$ws.Range("M136").Value = 392.2856999999999  # was 15
